$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Insert a new row at position 12 (pushes the old rows 12..33 down to 13..34),
# matching the new "dayofyear" CALC variable added to the data model.
$ws.Rows.Item(12).Insert() | Out-Null

$ws.Cells.Item(12, 1).Value = "CALC"
$ws.Cells.Item(12, 2).Value = "ready"
$ws.Cells.Item(12, 3).Value = 1
$ws.Cells.Item(12, 4).Value = "dayofyear"
$ws.Cells.Item(12, 5).Value = "metr"

# Re-apply the AutoFilter so its range grows to cover the new row (A1:G35),
# matching the worksheet's extra trailing blank filter row.
$ws.AutoFilterMode = $false
$ws.Range("A1:G35").AutoFilter() | Out-Null

# Keep the _FilterDatabase defined name in sync with the expanded filter range.
foreach ($n in $wb.Names) {
  if ($n.Name -eq "Tabelle1!_FilterDatabase") {
    $n.RefersTo = "=Tabelle1!`$A`$1:`$G`$35"
  }
}

# Move the active selection, mirroring the author's cursor position after the edit.
$ws.Range("E10").Select() | Out-Null
